$d = $word.ActiveDocument
$word.UserName = "Handley, Mark"
$d.TrackRevisions = $false

$r = $d.Range(0,0)
$r.InsertBefore("Some secret text from an old edit.`r")

$r0 = $d.Range(0,0)
$r0.InsertBefore("Z")
$d.Bookmarks.Add("_GoBack", $d.Range(0,1))
$d.Range(0,1).Delete()

$d.TrackRevisions = $true
$delRange = $d.Paragraphs(1).Range
$delRange.Delete()
